# Final Version of Code before Submission
# Update computed results (stress/strain cascade from a re-run of the analysis)
# and relabel "Local" -> "Global" Stress/Strain per final submission.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel section headers -------------------------------------------------
$ws.Range("A13").Value = "Global Stress"
$ws.Range("A18").Value = "Global Strain"

# --- Updated numeric results --------------------------------------------------
$newValues = @{
    "B1" = "829.65806826558526"
    "B2" = "-170.34193173441474"
    "B3" = "2.4980018054066022e-14"
    "B6" = "-0.02"
    "C6" = "-0.014999999999999999"
    "D6" = "-0.01"
    "E6" = "-0.0050000000000000001"
    "G6" = "0.0050000000000000001"
    "H6" = "0.01"
    "I6" = "0.014999999999999999"
    "J6" = "0.02"
    "B9" = "0.0031298277647945406"
    "E9" = "9.774755977538807e-35"
    "B10" = "-0.00090878954396307046"
    "E10" = "5.4228149378825262e-35"
    "B11" = "-2.6105912990350568e-19"
    "E11" = "-6.3888867485851463e-34"
    "B14" = "16538.617308757606"
    "C14" = "16538.617308757606"
    "D14" = "13680.754508811911"
    "E14" = "53242.010873672873"
    "F14" = "53242.010873672873"
    "G14" = "53242.010873672873"
    "H14" = "13680.754508811911"
    "I14" = "16538.617308757606"
    "J14" = "16538.617308757606"
    "B15" = "8461.3826912423829"
    "C15" = "8461.3826912423865"
    "D15" = "-26275.090891553962"
    "E15" = "9352.32550906918"
    "F15" = "9352.32550906918"
    "G15" = "9352.32550906918"
    "H15" = "-26275.090891553962"
    "I15" = "8461.3826912423865"
    "J15" = "8461.3826912423829"
    "B16" = "983.45999105945702"
    "C16" = "-983.45999105946066"
    "D16" = "2.9648969638288924e-14"
    "E16" = "2.9648969638288924e-14"
    "F16" = "2.9648969638288924e-14"
    "G16" = "2.9648969638288924e-14"
    "H16" = "2.9648969638288924e-14"
    "I16" = "-983.45999105946066"
    "J16" = "983.45999105945702"
    "B19" = "0.0060139918238021409"
    "C19" = "0.0060139918238021391"
    "D19" = "0.0095014918238021401"
    "E19" = "0.0025264918238021407"
    "F19" = "0.0025264918238021407"
    "G19" = "0.0025264918238021407"
    "H19" = "0.0095014918238021401"
    "I19" = "0.0060139918238021391"
    "J19" = "0.0060139918238021409"
    "B20" = "0.0019753745150445281"
    "C20" = "0.0019753745150445281"
    "D20" = "-0.0015121254849554708"
    "E20" = "0.0054628745150445291"
    "F20" = "0.0054628745150445291"
    "G20" = "0.0054628745150445291"
    "H20" = "-0.0015121254849554708"
    "I20" = "0.0019753745150445281"
    "J20" = "0.0019753745150445281"
    "B21" = "-0.0069749999999999986"
    "C21" = "0.0069749999999999977"
    "D21" = "2.9648969638288912e-20"
    "E21" = "2.9648969638288912e-20"
    "F21" = "2.9648969638288912e-20"
    "G21" = "2.9648969638288912e-20"
    "H21" = "2.9648969638288912e-20"
    "I21" = "0.0069749999999999977"
    "J21" = "-0.0069749999999999986"
}

foreach ($ref in $newValues.Keys) {
    $ws.Range($ref).Value = [double]$newValues[$ref]
}

# --- Column width tweaks (D, F, G narrowed slightly) --------------------------
$ws.Columns.Item(4).ColumnWidth = 14.7109375
$ws.Columns.Item(6).ColumnWidth = 14.7109375
$ws.Columns.Item(7).ColumnWidth = 14.7109375
